# Scheduled-runner market-price refresh for Sheets (Lamia_Profits workbook).
# Updates currentAveragePrice / crafting-cost / profit columns (H-N) per leve
# row across all eight job sheets to the latest Universalis price snapshot.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 4
$ws.Range("H4").Value = 255.76923
$ws.Range("I4").Value = 165.90909
$ws.Range("K4").Value = 165.90909
$ws.Range("M4").Value = -51.90908999999999
# Row 18
$ws.Range("H18").Value = 3799.3635
$ws.Range("I18").Value = 3799.3635
$ws.Range("K18").Value = 3799.3635
$ws.Range("M18").Value = -3515.3635
# Row 32
$ws.Range("H32").Value = 12844.6
$ws.Range("I32").Value = 14240.333
$ws.Range("J32").Value = 10751
$ws.Range("K32").Value = 14240.333
$ws.Range("L32").Value = 10751
$ws.Range("M32").Value = -13914.333
$ws.Range("N32").Value = -11403
# Row 112
$ws.Range("H112").Value = 1290.24
$ws.Range("J112").Value = 1398.4286
$ws.Range("L112").Value = 4195.2858
$ws.Range("N112").Value = -6411.2858
# Row 137
$ws.Range("H137").Value = 4959.8237
$ws.Range("J137").Value = 5208.5454
$ws.Range("L137").Value = 15625.6362
$ws.Range("N137").Value = -20725.6362
# Row 138
$ws.Range("H138").Value = 2845.894
$ws.Range("I138").Value = 1447.7273
$ws.Range("K138").Value = 4343.1819
$ws.Range("M138").Value = 796.8181000000004

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 44
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
# Row 97
$ws.Range("H97").Value = 899.2174
$ws.Range("I97").Value = 799.05554
$ws.Range("K97").Value = 799.05554
$ws.Range("M97").Value = -303.05554
# Row 132
$ws.Range("H132").Value = 2141.1072
$ws.Range("I132").Value = 1567.9615
$ws.Range("J132").Value = 9592
$ws.Range("K132").Value = 4703.8845
$ws.Range("L132").Value = 28776
$ws.Range("M132").Value = -2173.8845
$ws.Range("N132").Value = -33836

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 4017.111
$ws.Range("I86").Value = 3391.3333
$ws.Range("K86").Value = 3391.3333
$ws.Range("M86").Value = -2268.3333
# Row 89
$ws.Range("H89").Value = 4017.111
$ws.Range("I89").Value = 3391.3333
$ws.Range("K89").Value = 16956.6665
$ws.Range("M89").Value = -11340.6665
# Row 94
$ws.Range("H94").Value = 1860.5294
$ws.Range("I94").Value = 1330.6428
$ws.Range("K94").Value = 1330.6428
$ws.Range("M94").Value = -879.6428000000001
# Row 134
$ws.Range("H134").Value = 1545.1765
$ws.Range("I134").Value = 1317.8667
$ws.Range("K134").Value = 3953.6001
$ws.Range("M134").Value = -1418.6001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 2843.2222
$ws.Range("I22").Value = 298.57144
$ws.Range("K22").Value = 298.57144
$ws.Range("M22").Value = 51.42856
# Row 31
$ws.Range("H31").Value = 24880.432
$ws.Range("I31").Value = 3050.389
$ws.Range("K31").Value = 3050.389
$ws.Range("M31").Value = -2755.389
# Row 34
$ws.Range("H34").Value = 24880.432
$ws.Range("I34").Value = 3050.389
$ws.Range("K34").Value = 3050.389
$ws.Range("M34").Value = -2848.389
# Row 105
$ws.Range("H105").Value = 1721.8
$ws.Range("I105").Value = 874.7222
$ws.Range("J105").Value = 2992.4167
$ws.Range("K105").Value = 874.7222
$ws.Range("L105").Value = 2992.4167
$ws.Range("M105").Value = 872.2778
$ws.Range("N105").Value = -6486.4167
# Row 134
$ws.Range("H134").Value = 2739.7188
$ws.Range("I134").Value = 2088.5667
$ws.Range("K134").Value = 6265.7001
$ws.Range("M134").Value = -3730.7001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 7
$ws.Range("H7").Value = 229
$ws.Range("J7").Value = 130
$ws.Range("L7").Value = 390
$ws.Range("N7").Value = -614
# Row 34
$ws.Range("H34").Value = 1964.9333
$ws.Range("I34").Value = 1619.3334
$ws.Range("J34").Value = 2483.3333
$ws.Range("K34").Value = 4858.0002
$ws.Range("L34").Value = 7449.999899999999
$ws.Range("M34").Value = -4774.0002
$ws.Range("N34").Value = -7617.999899999999
# Row 46
$ws.Range("H46").Value = 2820
$ws.Range("I46").Value = 95
$ws.Range("K46").Value = 285
$ws.Range("M46").Value = -194
# Row 55
$ws.Range("H55").Value = 1971.2142
$ws.Range("I55").Value = 1409
$ws.Range("K55").Value = 4227
$ws.Range("M55").Value = -4050
# Row 131
$ws.Range("H131").Value = 6144981.5
$ws.Range("J131").Value = 4631501.5
$ws.Range("L131").Value = 13894504.5
$ws.Range("N131").Value = -13904584.5
# Row 137
$ws.Range("H137").Value = 85404.164
$ws.Range("J137").Value = 127250
$ws.Range("L137").Value = 381750
$ws.Range("N137").Value = -391950
# Row 140
$ws.Range("H140").Value = 2560.923
$ws.Range("I140").Value = 2224.8333
$ws.Range("K140").Value = 6674.499899999999
$ws.Range("M140").Value = -1494.499899999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 62
$ws.Range("H62").Value = 37781.332
$ws.Range("I62").Value = 37781.332
$ws.Range("K62").Value = 37781.332
$ws.Range("M62").Value = -37095.332
# Row 65
$ws.Range("H65").Value = 37781.332
$ws.Range("I65").Value = 37781.332
$ws.Range("K65").Value = 113343.996
$ws.Range("M65").Value = -109911.996
# Row 80
$ws.Range("H80").Value = 231072.36
$ws.Range("I80").Value = 456709
$ws.Range("K80").Value = 456709
$ws.Range("M80").Value = -455711
# Row 83
$ws.Range("H83").Value = 231072.36
$ws.Range("I83").Value = 456709
$ws.Range("K83").Value = 2283545
$ws.Range("M83").Value = -2278553
# Row 132
$ws.Range("H132").Value = 5958.3945
$ws.Range("I132").Value = 5233.472
$ws.Range("K132").Value = 15700.416
$ws.Range("M132").Value = -13170.416

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 25
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("M25").ClearContents()
$ws.Range("N25").ClearContents()
# Row 46
$ws.Range("H46").Value = 3425.963
$ws.Range("I46").Value = 2077.889
$ws.Range("J46").Value = 4100
$ws.Range("K46").Value = 2077.889
$ws.Range("L46").Value = 4100
$ws.Range("M46").Value = -1889.889
$ws.Range("N46").Value = -4476
# Row 48
$ws.Range("H48").Value = 29999
$ws.Range("J48").Value = 29999
$ws.Range("L48").Value = 29999
$ws.Range("N48").Value = -31321
# Row 55
$ws.Range("H55").Value = 2175967.5
$ws.Range("J55").Value = 4823.25
$ws.Range("L55").Value = 4823.25
$ws.Range("N55").Value = -5169.25
# Row 68
$ws.Range("H68").Value = 4689.5454
$ws.Range("I68").Value = 3766.3333
$ws.Range("J68").Value = 5797.4
$ws.Range("K68").Value = 3766.3333
$ws.Range("L68").Value = 5797.4
$ws.Range("M68").Value = -3017.3333
$ws.Range("N68").Value = -7295.4
# Row 71
$ws.Range("H71").Value = 4689.5454
$ws.Range("I71").Value = 3766.3333
$ws.Range("J71").Value = 5797.4
$ws.Range("K71").Value = 18831.6665
$ws.Range("L71").Value = 28987
$ws.Range("M71").Value = -15087.6665
$ws.Range("N71").Value = -36475
# Row 132
$ws.Range("H132").Value = 7621.5713
$ws.Range("I132").Value = 6869.095
$ws.Range("K132").Value = 20607.285
$ws.Range("M132").Value = -18077.285
# Row 136
$ws.Range("H136").Value = 6859.6206
$ws.Range("I136").Value = 2207.111
$ws.Range("J136").Value = 14472.818
$ws.Range("K136").Value = 6621.333
$ws.Range("L136").Value = 43418.454
$ws.Range("M136").Value = -4071.333
$ws.Range("N136").Value = -48518.454

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 26
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()
# Row 107
$ws.Range("H107").Value = 940.3333
$ws.Range("I107").Value = 745.375
$ws.Range("K107").Value = 2236.125
$ws.Range("M107").Value = -316.125
# Row 132
$ws.Range("H132").Value = 3163.6191
$ws.Range("I132").Value = 2158.1765
$ws.Range("K132").Value = 6474.529500000001
$ws.Range("M132").Value = -3944.529500000001
# Row 136
$ws.Range("H136").Value = 4058.7222
$ws.Range("I136").Value = 3190.75
$ws.Range("K136").Value = 9572.25
$ws.Range("M136").Value = -7022.25
